# Weekly driver report update for 2025-04-20
# Updates the "Bad Drivers" summary row and re-ranks / refreshes the
# "Good Drivers" table (adapter name, client/sample counts, roaming %,
# and driver-vintage date) on the "Driver Summary" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Bad Drivers summary (rows 3-4) ---------------------------------
$ws.Range("C3").Value = 164
$ws.Range("D3").Value = 98.8
$ws.Range("C4").Value = 164

# --- Helper: write a "Driver Vintage" date as literal text ----------
# A plain `.Value = "yyyy-mm-dd"` assignment gets auto-recognized as a
# real date (and reformatted), so the target cell is first marked as
# Text, written, and then has its number format restored by pasting
# the format from its same-row "Good Roaming Calculation" neighbor
# (column D), which already carries the plain style used throughout
# this table.
function Set-VintageText($cell, $text) {
    $col = $ws.Range($cell)
    $col.NumberFormat = "@"
    $col.Value = $text
    $fmtSource = $ws.Range("D" + $cell.Substring(1))
    $fmtSource.Copy()
    $col.PasteSpecial(-4122)
    $excel.CutCopyMode = $false
}

# --- Good Drivers table (rows 12-17) ---------------------------------
$ws.Range("A12").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.100.0.4"
$ws.Range("B12").Value = 445055
$ws.Range("D12").Value = 99.90000000000001
Set-VintageText "E12" "2024-11-10"

$ws.Range("A13").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.80.0.9"
$ws.Range("B13").Value = 77849
$ws.Range("D13").Value = 99.90000000000001
Set-VintageText "E13" "2021-08-18"

$ws.Range("A14").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.50.1.1"
$ws.Range("B14").Value = 34244
$ws.Range("D14").Value = 100
Set-VintageText "E14" "2021-04-27"

$ws.Range("A15").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.110.3.2"
$ws.Range("B15").Value = 59673
$ws.Range("D15").Value = 100
Set-VintageText "E15" "2020-08-05"

$ws.Range("A16").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.70.0.6"
$ws.Range("B16").Value = 113652
$ws.Range("D16").Value = 100
Set-VintageText "E16" "2020-01-06"

$ws.Range("A17").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.60.2.1"
$ws.Range("B17").Value = 56018
$ws.Range("D17").Value = 100
# E17 stays "2019-12-14" - unchanged by this week's refresh.
